$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Mark existing rows 44-47 with new status / extra columns ----

# Row 44: TitleScreen - fix mesh rendering -> DONE
$ws.Range("D44").Value2 = "DONE"

# Row 45: TitleScreen - fix skybox rendering -> DONE
$ws.Range("D45").Value2 = "DONE"

# Row 46: TitleScreen - fix placement of mesh and skybox choice -> DONE
$ws.Range("D46").Value2 = "DONE"
$ws.Range("E46").Value2 = "TODO: update controls/have controls blinking on bottom of screen or something"
$ws.Range("F46").Value2 = "DONE"

# Row 47: TItleScreen - fix camera -> In progress / Fiddling
$ws.Range("D47").Value2 = "In progress"
$ws.Range("E47").Value2 = "Fiddling"

# ---- New rows 48-51 ----

# Row 48
$ws.Range("B48").Value2 = "CarScreen - fix the ""Player X…. "" text so that doesn't switch when done"
$ws.Range("C48").Value2 = 42941
$ws.Range("C48").NumberFormat = "d-mmm-yy"

# Row 49
$ws.Range("B49").Value2 = "TitleScreen - ""PLAY"" button before going to start of game?"
$ws.Range("C49").Value2 = 42941
$ws.Range("C49").NumberFormat = "d-mmm-yy"

# Row 50
$ws.Range("B50").Value2 = "TitleScreen - fix colour of text, black sort of hard to see"
$ws.Range("C50").Value2 = 42941
$ws.Range("C50").NumberFormat = "d-mmm-yy"

# Row 51
$ws.Range("B51").Value2 = "Fix pause menu"

# ---- Update selection to match final cursor position ----
$ws.Range("B51").Select()
